# Auto-generated Excel COM-interop script to apply updated profit calculations
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 602.7273
$ws.Range("I53").Value = 162.16667
$ws.Range("J53").Value = 1131.4
$ws.Range("K53").Value = 162.16667
$ws.Range("L53").Value = 1131.4
$ws.Range("M53").Value = 474.83333
$ws.Range("N53").Value = -2405.4
$ws.Range("H64").Value = 42338.46
$ws.Range("I64").Value = 4200
$ws.Range("J64").Value = 86833.336
$ws.Range("K64").Value = 4200
$ws.Range("L64").Value = 86833.336
$ws.Range("M64").Value = -3952
$ws.Range("N64").Value = -87329.336
$ws.Range("H67").Value = 42338.46
$ws.Range("I67").Value = 4200
$ws.Range("J67").Value = 86833.336
$ws.Range("K67").Value = 4200
$ws.Range("L67").Value = 86833.336
$ws.Range("M67").Value = -3342
$ws.Range("N67").Value = -88549.336
$ws.Range("H113").Value = 1925
$ws.Range("I113").Value = 2350
$ws.Range("J113").Value = 1742.8572
$ws.Range("K113").Value = 2350
$ws.Range("L113").Value = 1742.8572
$ws.Range("M113").Value = 904
$ws.Range("N113").Value = -8250.8572
$ws.Range("H121").Value = 2499.5
$ws.Range("J121").Value = 2499.5
$ws.Range("L121").Value = 7498.5
$ws.Range("N121").Value = -10992.5
$ws.Range("H136").Value = 54890
$ws.Range("J136").Value = 54890
$ws.Range("L136").Value = 54890
$ws.Range("N136").Value = -65090
$ws.Range("H141").Value = 1408.2307
$ws.Range("I141").Value = 1087.2222
$ws.Range("J141").Value = 2130.5
$ws.Range("K141").Value = 3261.6666
$ws.Range("L141").Value = 6391.5
$ws.Range("M141").Value = 1918.3334
$ws.Range("N141").Value = -16751.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 27457.715
$ws.Range("J7").Value = 27457.715
$ws.Range("L7").Value = 27457.715
$ws.Range("N7").Value = -27685.715
$ws.Range("H32").Value = 6921.707
$ws.Range("I32").Value = 4012.8635
$ws.Range("K32").Value = 4012.8635
$ws.Range("M32").Value = -3725.8635
$ws.Range("H45").Value = 806.8461
$ws.Range("I45").Value = 721
$ws.Range("K45").Value = 721
$ws.Range("M45").Value = -344
$ws.Range("H110").Value = 1088.24
$ws.Range("I110").Value = 998.5625
$ws.Range("J110").Value = 1247.6666
$ws.Range("K110").Value = 998.5625
$ws.Range("L110").Value = 1247.6666
$ws.Range("M110").Value = 1046.4375
$ws.Range("N110").Value = -5337.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1916.7858
$ws.Range("I20").Value = 1711.7391
$ws.Range("K20").Value = 1711.7391
$ws.Range("M20").Value = -1464.7391
$ws.Range("H86").Value = 2382.9714
$ws.Range("I86").Value = 2345.5
$ws.Range("J86").Value = 2446.3845
$ws.Range("K86").Value = 2345.5
$ws.Range("L86").Value = 2446.3845
$ws.Range("M86").Value = -1222.5
$ws.Range("N86").Value = -4692.3845
$ws.Range("H89").Value = 2382.9714
$ws.Range("I89").Value = 2345.5
$ws.Range("J89").Value = 2446.3845
$ws.Range("K89").Value = 11727.5
$ws.Range("L89").Value = 12231.9225
$ws.Range("M89").Value = -6111.5
$ws.Range("N89").Value = -23463.9225
$ws.Range("H98").Value = 59385.5
$ws.Range("J98").Value = 59385.5
$ws.Range("L98").Value = 59385.5
$ws.Range("N98").Value = -65375.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1980.7273
$ws.Range("I58").Value = 1605.125
$ws.Range("K58").Value = 1605.125
$ws.Range("M58").Value = -1402.125
$ws.Range("H86").Value = 2769.8076
$ws.Range("I86").Value = 2570.6924
$ws.Range("J86").Value = 2968.923
$ws.Range("K86").Value = 2570.6924
$ws.Range("L86").Value = 2968.923
$ws.Range("M86").Value = -1447.6924
$ws.Range("N86").Value = -5214.923
$ws.Range("H89").Value = 2769.8076
$ws.Range("I89").Value = 2570.6924
$ws.Range("J89").Value = 2968.923
$ws.Range("K89").Value = 12853.462
$ws.Range("L89").Value = 14844.615
$ws.Range("M89").Value = -7237.462
$ws.Range("N89").Value = -26076.615
$ws.Range("H127").Value = 59773.332
$ws.Range("J127").Value = 59773.332
$ws.Range("L127").Value = 59773.332
$ws.Range("N127").Value = -69693.33199999999
$ws.Range("H129").Value = 53947.6
$ws.Range("J129").Value = 53947.6
$ws.Range("L129").Value = 53947.6
$ws.Range("N129").Value = -63947.6
$ws.Range("H136").Value = 1980.7273
$ws.Range("I136").Value = 1605.125
$ws.Range("K136").Value = 4815.375
$ws.Range("M136").Value = -2265.375
$ws.Range("H139").Value = 21342.615
$ws.Range("J139").Value = 21342.615
$ws.Range("L139").Value = 21342.615
$ws.Range("N139").Value = -31622.615

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 60.875
$ws.Range("I23").Value = 72.5
$ws.Range("J23").Value = 53.9
$ws.Range("K23").Value = 217.5
$ws.Range("L23").Value = 161.7
$ws.Range("M23").Value = 17.5
$ws.Range("N23").Value = -631.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 29498.334
$ws.Range("J54").Value = 29498.334
$ws.Range("L54").Value = 29498.334
$ws.Range("N54").Value = -30278.334
$ws.Range("H70").Value = 4840.4585
$ws.Range("I70").Value = 4938.9443
$ws.Range("J70").Value = 4545
$ws.Range("K70").Value = 4938.9443
$ws.Range("L70").Value = 4545
$ws.Range("M70").Value = -4668.9443
$ws.Range("N70").Value = -5085
$ws.Range("H73").Value = 4840.4585
$ws.Range("I73").Value = 4938.9443
$ws.Range("J73").Value = 4545
$ws.Range("K73").Value = 4938.9443
$ws.Range("L73").Value = 4545
$ws.Range("M73").Value = -4002.9443
$ws.Range("N73").Value = -6417

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 616.6667
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -255
$ws.Range("N22").Value = -1340
$ws.Range("H27").Value = 616.6667
$ws.Range("I27").Value = 550
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 550
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -443
$ws.Range("N27").Value = -964
$ws.Range("H93").Value = 1468.7273
$ws.Range("I93").Value = 1584.3334
$ws.Range("J93").Value = 1330
$ws.Range("K93").Value = 1584.3334
$ws.Range("L93").Value = 1330
$ws.Range("M93").Value = -336.3334
$ws.Range("N93").Value = -3826
$ws.Range("H122").Value = 9196.809999999999
$ws.Range("I122").Value = 13598
$ws.Range("J122").Value = 5195.727
$ws.Range("K122").Value = 40794
$ws.Range("L122").Value = 15587.181
$ws.Range("M122").Value = -38344
$ws.Range("N122").Value = -20487.181
$ws.Range("H136").Value = 1610.56
$ws.Range("I136").Value = 769.1177
$ws.Range("J136").Value = 3398.625
$ws.Range("K136").Value = 2307.3531
$ws.Range("L136").Value = 10195.875
$ws.Range("M136").Value = 242.6468999999997
$ws.Range("N136").Value = -15295.875
$ws.Range("H141").Value = 49250
$ws.Range("J141").Value = 48500
$ws.Range("L141").Value = 48500
$ws.Range("N141").Value = -58860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2853.8333
$ws.Range("I81").Value = 1035.3636
$ws.Range("K81").Value = 2070.7272
$ws.Range("M81").Value = -1009.7272
$ws.Range("H84").Value = 2853.8333
$ws.Range("I84").Value = 1035.3636
$ws.Range("K84").Value = 10353.636
$ws.Range("M84").Value = -5049.635999999999
